# Add a new "Franck" row (idPers=17) to the "Personnes" sheet,
# mirroring the existing rows (A: id number, B: name string, C: left blank).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Franck"
$ws.Range("C18").NumberFormat = "General"
